$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose column D value ("Equipo Programático") should become the new,
# shorter label "E. Programático" (a new shared string entry).
$rows = @(11,12,13,14,15,16,17,18,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,65,66,67,100,101)

foreach ($r in $rows) {
    $ws.Range("D$r").Value = "E. Programático"
}

# Update the active selection shown in the sheet view.
[void]$ws.Range("D20").Select()
